# Build site at 2023-04-12 14:53:07 UTC
# Reworks the "Objetivos/Docentes/Programa resumido/Programa/Bibliografia"
# block of the LOQ4037 course sheet: fixes several rows whose B/C values had
# been shifted/duplicated, and appends the Portuguese "Programa" syllabus
# body, the Portuguese "Bibliografia" references and a second "Requisitos"
# line (LOQ4098) that were previously missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Row heights that need to change.
#    AutoFit() fully clears any explicit/custom height (rows that must end
#    up with the sheet's default height again); RowHeight = N sets an
#    explicit custom height (matches ht="N" customHeight="1" in the xml).
# ---------------------------------------------------------------------
$ws.Rows.Item(13).AutoFit()      | Out-Null   # was ht=60 -> default
$ws.Rows.Item(15).RowHeight = 60             # was ht=120 -> 60
$ws.Rows.Item(17).RowHeight = 120            # was default -> 120
$ws.Rows.Item(18).AutoFit()      | Out-Null   # was ht=60 -> default
$ws.Rows.Item(21).RowHeight = 60             # was ht=120 -> 60
$ws.Rows.Item(22).RowHeight = 120            # was default -> 120
$ws.Rows.Item(23).AutoFit()      | Out-Null   # was ht=30 -> default
$ws.Rows.Item(25).RowHeight = 30             # new row -> 30

# ---------------------------------------------------------------------
# 2) Pre-format the brand-new cells by copying number/font/alignment
#    formats from stable donor cells that keep the same role throughout
#    (column A labels = bold/top; column B = top+wrap; column C = top+wrap
#    red) before any values are written into them.
# ---------------------------------------------------------------------
Set-CellFormat "B3" "B17"
Set-CellFormat "C3" "C17"
Set-CellFormat "B3" "B22"
Set-CellFormat "C3" "C22"
Set-CellFormat "A12" "A23"
Set-CellFormat "B3" "B25"
Set-CellFormat "C3" "C25"
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Remove the cells that no longer hold data in the new layout.
#    Range.Clear() drops the cell (value + style) entirely, matching the
#    target xml where these cells simply don't exist any more.
# ---------------------------------------------------------------------
$ws.Range("A13").Clear() | Out-Null
$ws.Range("B18").Clear() | Out-Null
$ws.Range("C18").Clear() | Out-Null
$ws.Range("B23").Clear() | Out-Null
$ws.Range("C23").Clear() | Out-Null

# ---------------------------------------------------------------------
# 4) Write the corrected / shifted text into place (existing cells whose
#    value changes, and the brand new cells formatted above).
# ---------------------------------------------------------------------

# Row 10 - "Objetivos:" body (was wrongly holding the "Docentes" name).
$objetivosPt = "Gerais - Apresentar e Ensinar conceitos de Química Orgânica como instrumentos importantes para a compreensão de estratégias e operações industriais e tecnológicas. Abordar problemáticas sociais e ambientais com as quais a engenharia química está relacionada, tornando-os dessa forma, aptos a exercerem a função de Engenheiro Químico, e realizarem as mudanças que se façam necessárias.Específicos – Compreender e descrever o mecanismo das reações orgânicas e a sua importância para o aprimoramento e desenvolvimento de processos industriais sintéticos e de etapas de formulação. Aprofundar o conceito de estrutura-reatividade e propriedades dos materiais."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# Row 13 - "Docentes responsáveis:" value (now under the A12 label).
$docente = "210064 - Eduardo Rezende Triboni"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# Row 14 - "Programa resumido:" (label shifted down, value replaced).
$ws.Range("A14").Value = "Programa resumido:"
$resumoPt = "Propriedade gerais dos compostos orgânicos. Estrutura, métodos de obtenção, propriedades físicas, reações dos hidrocarbonetos alifáticos e aromáticos, haletos orgânicos, álcoois e características estruturaiscomo estereoquímica e a relação estrutura-reatividade."
$ws.Range("B14").Value = $resumoPt
$ws.Range("C14").Value = $resumoPt

# Row 15 - "Short syllabus:" (label shifted down, value fixed - was a date).
$ws.Range("A15").Value = "Short syllabus:"
$shortSyllabus = "General property of organic compounds. Physical properties, reactions of aliphatic and aromatic hydrocarbons, organic halides, ethers, alcohols and structural characteristics as stereochemistry and structure-reactivity."
$ws.Range("B15").Value = $shortSyllabus
$ws.Range("C15").Value = $shortSyllabus

# Row 16 - "Programa:" (label shifted down, Portuguese syllabus body).
$ws.Range("A16").Value = "Programa:"
$programaPt = "1.Teoria de Bronsted e de Lewis e acidez de compostos orgânicos2.Alcanos - Processos de obtenção, Propriedades físicas, Análise Conformacional. Reação de Substituição Radicalar. 3.Isomeria Constitucional e Isomeria Espacial (Estereoquímica). Quiralidade, Nomenclatura R/S, classificação de estereoisômeros. Polarímetro e Técnicas de  Resolução de Isômeros Espaciais.4.Haletos de Alquila – Substituição Nucleofílica, SN1, SN2, E1, E2. 5.Alcenos, Alcadienos e Alcinos – Propriedades físicas e químicas. Reação de adição eletrofílica (hidroalogenação, Hidratação, Halogenação, Diels-Alder, Redução-Oxidação). Adição conjugada em dienos (produto termodinâmico e cinético) 6. Fundamentos de RMN, Infra-vermelho, Ultra-violeta e Fluorescencia 7.Compostos aromáticos – Propriedades físicas dos aromáticos. Reações de Substituição Eletrofílica Aromática. Efeito de Grupos Substituintes. Reação de Substituição Nucleofílica.8.Álcoois e Éteres – Propriedades físicas, reações e mecanismos."
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# Row 17 - "Syllabus:" (new row; English syllabus body moved here).
$ws.Range("A17").Value = "Syllabus:"
$syllabusEn = "1.Bronsted and Lewis acid of the organic compounds2.Alkanes - obtaining processes, physical properties, conformational analysis. Radical substitution reaction.3.Constitutional isomerism and Stereochemistry. Chirality, nomenclature R/S. Polarimeter and Techniques for resolution of stereoisomers.4.Alkyl halides - Nucleophilic Substitution, SN1, SN2, E1, E2.5.Alkenes, alkadienes and alkynes - Physical and chemical properties. Electrophilic addition reaction (hidrohalogenation, hydration, halogenation, Diels-Alder, reduction and oxidation). Conjugated Addition in dienes (thermodynamic and kinetic product).6 Background of NMR, InfraRed, UV and Fluorescence techniques7.Aromatic compounds - Physical properties. Aromatic Eletrophilic Substitution . Effect of Substituent Groups. Aromatic Nucleophilic Substitution.8.Alcohols and ethers - physical properties, reactions and mechanisms."
$ws.Range("B17").Value = $syllabusEn
$ws.Range("C17").Value = $syllabusEn

# Row 18 - "Avaliação:" (label shifted down, no B/C content any more).
$ws.Range("A18").Value = "Avaliação:"

# Row 19 - "Método:" label only shifts down, B/C keep their existing text.
$ws.Range("A19").Value = "Método:"

# Row 20 - "Critério:" label shifts down, B/C keep their existing text.
$ws.Range("A20").Value = "Critério:"

# Row 21 - "Norma de recuperação:" label shifts down, B/C keep existing text.
$ws.Range("A21").Value = "Norma de recuperação:"

# Row 22 - "Bibliografia:" (new label here) + reference list (new content).
$ws.Range("A22").Value = "Bibliografia:"
$bibliografia = @"
BRESLOW, R. Questões e Exercícios de Química Orgânica. São Paulo: Makrons Books Editora, 1996. 
BRUICE, P. Y. Química Orgânica, vol 1 e 2, São Paulo: Editora Pearson Prentice Hall, 2006. 
HENDRIKSON, James B.; CRAM, Donald J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966.
MCMURRY, John. Química Orgânica. São Paulo: Editora Pioneira Thomson Leraning, 2005.
SOLOMONS, T.W.G; FRYHLE, Graig. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001.
MORRISON, R.; BOYD, R. Química Orgânica. São Paulo: Editora Calouste Gulbenkian, 2008.
"@
$bibliografia = $bibliografia.TrimEnd("`r","`n")
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# Row 23 - "Requisitos:" label (new A cell; no B/C on this row any more).
$ws.Range("A23").Value = "Requisitos:"

# Row 24 - keeps the LOQ4097 requirement text that used to sit in row 23.
$req1 = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`r`n"
$ws.Range("B24").Value = $req1
$ws.Range("C24").Value = $req1

# Row 25 - brand new row holding the LOQ4098 requirement text.
$req2 = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`r`n"
$ws.Range("B25").Value = $req2
$ws.Range("C25").Value = $req2

# ---------------------------------------------------------------------
# 5) Column layout: column A's custom width used to bleed into column B
#    (min=1 max=2); narrow that range down to column A only so column B
#    falls back to its own (wider) width rule.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.83
$ws.Columns.Item(2).ColumnWidth = 59.83
